$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price data rows got reshuffled (dates / volumes / prices / unit
# descriptions / $/Kg moved between rows). Row 6 is untouched.
# Apply the new values per row, for columns D, M, N, O, P, Q, S.
# Single-quoted strings are used so the literal '$' in the unit text is not
# treated as PowerShell variable interpolation.

# Row 2 <- old Row 5 values
$ws.Range("D2").Value = 44162
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("Q2").Value = '$/caja 14 kilos empedrada'
$ws.Range("S2").Value = 500

# Row 3 <- old Row 7 values
$ws.Range("D3").Value = 44176
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("Q3").Value = '$/caja 14 kilos empedrada'
$ws.Range("S3").Value = 500

# Row 4 <- old Row 8 values
$ws.Range("D4").Value = 44400
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("Q4").Value = '$/caja 14 kilos'
$ws.Range("S4").Value = 714

# Row 5 <- old Row 2 values
$ws.Range("D5").Value = 44309
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("Q5").Value = '$/caja 14 kilos empedrada'
$ws.Range("S5").Value = 500

# Row 7 <- old Row 9 values
$ws.Range("D7").Value = 44491
$ws.Range("M7").Value = 180
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 9000
$ws.Range("P7").Value = 9000
$ws.Range("Q7").Value = '$/caja 14 kilos empedrada'
$ws.Range("S7").Value = 643

# Row 8 <- old Row 4 values
$ws.Range("D8").Value = 44351
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("Q8").Value = '$/caja 14 kilos empedrada'
$ws.Range("S8").Value = 714

# Row 9 <- old Row 3 values
$ws.Range("D9").Value = 44397
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 11000
$ws.Range("O9").Value = 11000
$ws.Range("P9").Value = 11000
$ws.Range("Q9").Value = '$/caja 14 kilos'
$ws.Range("S9").Value = 786
